# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) per leve row across the Goblin_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6
$ws.Range("H6").Value = 8481.25
$ws.Range("I6").Value = 11142
$ws.Range("K6").Value = 33426
$ws.Range("M6").Value = -33314

# row 12
$ws.Range("H12").Value = 1999.5
$ws.Range("I12").Value = 1999
$ws.Range("K12").Value = 1999
$ws.Range("M12").Value = -1829

# row 15
$ws.Range("H15").Value = 2350.4814
$ws.Range("I15").Value = 2350.4814
$ws.Range("K15").Value = 7051.4442
$ws.Range("M15").Value = -6882.4442

# row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# row 58
$ws.Range("H58").Value = 9780.066000000001
$ws.Range("J58").Value = 48332.668
$ws.Range("L58").Value = 144998.004
$ws.Range("N58").Value = -145298.004

# row 80
$ws.Range("H80").Value = 1679.9524
$ws.Range("I80").Value = 755.8182
$ws.Range("K80").Value = 2267.4546
$ws.Range("M80").Value = -1269.4546

# row 83
$ws.Range("H83").Value = 1679.9524
$ws.Range("I83").Value = 755.8182
$ws.Range("K83").Value = 6802.3638
$ws.Range("M83").Value = -1810.3638

# row 87
$ws.Range("H87").Value = 157939400
$ws.Range("J87").Value = 157939400
$ws.Range("L87").Value = 157939400
$ws.Range("N87").Value = -157941896

# row 90
$ws.Range("H90").Value = 157939400
$ws.Range("J90").Value = 157939400
$ws.Range("L90").Value = 473818200
$ws.Range("N90").Value = -473830680

# row 92
$ws.Range("H92").Value = 1292
$ws.Range("I92").Value = 1168.9131
$ws.Range("K92").Value = 1168.9131
$ws.Range("M92").Value = 79.08690000000001

# row 100
$ws.Range("H100").Value = 5949.25
$ws.Range("J100").Value = 8000
$ws.Range("L100").Value = 8000
$ws.Range("N100").Value = -9082

# row 106
$ws.Range("H106").Value = 2586.8262
$ws.Range("I106").Value = 2293.0527
$ws.Range("K106").Value = 2293.0527
$ws.Range("M106").Value = -1662.0527

# row 133
$ws.Range("H133").Value = 101960.164
$ws.Range("J133").Value = 101960.164
$ws.Range("L133").Value = 101960.164
$ws.Range("N133").Value = -112080.164

# row 136
$ws.Range("H136").Value = 79873.664
$ws.Range("J136").Value = 79873.664
$ws.Range("L136").Value = 79873.664
$ws.Range("N136").Value = -90073.664

$ws = $wb.Worksheets.Item("ARM")
# row 63
$ws.Range("H63").Value = 12000
$ws.Range("J63").Value = 11800
$ws.Range("L63").Value = 11800
$ws.Range("N63").Value = -13172

# row 66
$ws.Range("H66").Value = 12000
$ws.Range("J66").Value = 11800
$ws.Range("L66").Value = 59000
$ws.Range("N66").Value = -65864

# row 88
$ws.Range("H88").Value = 9093687
$ws.Range("I88").Value = 25002612
$ws.Range("J88").Value = 2873.3572
$ws.Range("K88").Value = 25002612
$ws.Range("L88").Value = 2873.3572
$ws.Range("M88").Value = -25002206
$ws.Range("N88").Value = -3685.3572

# row 91
$ws.Range("H91").Value = 9093687
$ws.Range("I91").Value = 25002612
$ws.Range("J91").Value = 2873.3572
$ws.Range("K91").Value = 25002612
$ws.Range("L91").Value = 2873.3572
$ws.Range("M91").Value = -25001208
$ws.Range("N91").Value = -5681.3572

# row 92
$ws.Range("H92").Value = 45425
$ws.Range("J92").Value = 45425
$ws.Range("L92").Value = 45425
$ws.Range("N92").Value = -50417

# row 120
$ws.Range("H120").Value = 55000
$ws.Range("J120").Value = 55000
$ws.Range("L120").Value = 55000
$ws.Range("N120").Value = -64676

# row 133
$ws.Range("H133").Value = 74997
$ws.Range("J133").Value = 74997
$ws.Range("L133").Value = 74997
$ws.Range("N133").Value = -80057

$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 6895.08
$ws.Range("J86").Value = 2969
$ws.Range("L86").Value = 2969
$ws.Range("N86").Value = -5215

# row 89
$ws.Range("H89").Value = 6895.08
$ws.Range("J89").Value = 2969
$ws.Range("L89").Value = 14845
$ws.Range("N89").Value = -26077

# row 92
$ws.Range("H92").Value = 36940.1
$ws.Range("J92").Value = 36940.1
$ws.Range("L92").Value = 36940.1
$ws.Range("N92").Value = -41932.1

# row 94
$ws.Range("H94").Value = 2451.5
$ws.Range("I94").Value = 2662.1
$ws.Range("J94").Value = 1398.5
$ws.Range("K94").Value = 2662.1
$ws.Range("L94").Value = 1398.5
$ws.Range("M94").Value = -2211.1
$ws.Range("N94").Value = -2300.5

# row 132
$ws.Range("H132").Value = 213960.17
$ws.Range("J132").Value = 213960.17
$ws.Range("L132").Value = 213960.17
$ws.Range("N132").Value = -224080.17

# row 133
$ws.Range("H133").Value = 81933.336
$ws.Range("J133").Value = 81933.336
$ws.Range("L133").Value = 81933.336
$ws.Range("N133").Value = -92053.336

$ws = $wb.Worksheets.Item("CRP")
# row 17
$ws.Range("H17").Value = 7333.2856
$ws.Range("J17").Value = 7866.6
$ws.Range("L17").Value = 7866.6
$ws.Range("N17").Value = -8214.6

# row 25
$ws.Range("H25").Value = 2562.5
$ws.Range("I25").Value = 820
$ws.Range("K25").Value = 820
$ws.Range("M25").Value = -646

# row 59
$ws.Range("H59").Value = 25000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# row 60
$ws.Range("H60").Value = 3499.1667
$ws.Range("I60").Value = 3499.1667
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 3499.1667
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -2988.1667
$ws.Range("N60").ClearContents()

# row 122
$ws.Range("H122").Value = 3892.6843
$ws.Range("I122").Value = 3966.3125
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 11898.9375
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -9448.9375
$ws.Range("N122").Value = -15400

# row 141
$ws.Range("H141").Value = 292098.78
$ws.Range("J141").Value = 292098.78
$ws.Range("L141").Value = 292098.78
$ws.Range("N141").Value = -302458.78

$ws = $wb.Worksheets.Item("CUL")
# row 80
$ws.Range("H80").Value = 2339.7144
$ws.Range("J80").Value = 2339.7144
$ws.Range("L80").Value = 7019.1432
$ws.Range("N80").Value = -8891.143199999999

# row 83
$ws.Range("H83").Value = 2339.7144
$ws.Range("J83").Value = 2339.7144
$ws.Range("L83").Value = 21057.4296
$ws.Range("N83").Value = -30417.4296

$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 35717452
$ws.Range("I102").Value = 50001532
$ws.Range("K102").Value = 50001532
$ws.Range("M102").Value = -49999910

# row 104
$ws.Range("H104").Value = 43671
$ws.Range("J104").Value = 43671
$ws.Range("L104").Value = 43671
$ws.Range("N104").Value = -50659

# row 105
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988

# row 132
$ws.Range("H132").Value = 3120.2632
$ws.Range("I132").Value = 3074.75
$ws.Range("J132").Value = 3141.2693
$ws.Range("K132").Value = 9224.25
$ws.Range("L132").Value = 9423.8079
$ws.Range("M132").Value = -6694.25
$ws.Range("N132").Value = -14483.8079

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 816.6667
$ws.Range("I22").Value = 816.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 816.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -521.6667
$ws.Range("N22").ClearContents()

# row 27
$ws.Range("H27").Value = 816.6667
$ws.Range("I27").Value = 816.6667
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 816.6667
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -709.6667
$ws.Range("N27").ClearContents()

# row 120
$ws.Range("H120").Value = 101998
$ws.Range("J120").Value = 101998
$ws.Range("L120").Value = 101998
$ws.Range("N120").Value = -111674

# row 122
$ws.Range("H122").Value = 6303.2
$ws.Range("I122").Value = 6473.3184
$ws.Range("K122").Value = 19419.9552
$ws.Range("M122").Value = -16969.9552

# row 132
$ws.Range("H132").Value = 4008680.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4008680.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12026041.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -12031101.5

$ws = $wb.Worksheets.Item("WVR")
# row 6
$ws.Range("H6").Value = 12000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 12000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 12000
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -12230

# row 23
$ws.Range("H23").Value = 802.5
$ws.Range("I23").Value = 802.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 802.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -573.5
$ws.Range("N23").ClearContents()

# row 81
$ws.Range("H81").Value = 4966.5557
$ws.Range("I81").Value = 3483.3333
$ws.Range("J81").Value = 7933
$ws.Range("K81").Value = 6966.6666
$ws.Range("L81").Value = 15866
$ws.Range("M81").Value = -5905.6666
$ws.Range("N81").Value = -17988

# row 84
$ws.Range("H84").Value = 4966.5557
$ws.Range("I84").Value = 3483.3333
$ws.Range("J84").Value = 7933
$ws.Range("K84").Value = 34833.333
$ws.Range("L84").Value = 79330
$ws.Range("M84").Value = -29529.333
$ws.Range("N84").Value = -89938

# row 105
$ws.Range("H105").Value = 30333.334
$ws.Range("J105").Value = 30333.334
$ws.Range("L105").Value = 30333.334
$ws.Range("N105").Value = -37321.334

# row 122
$ws.Range("H122").Value = 3615.4055
$ws.Range("I122").Value = 1723.85
$ws.Range("J122").Value = 5840.7646
$ws.Range("K122").Value = 5171.549999999999
$ws.Range("L122").Value = 17522.2938
$ws.Range("M122").Value = -2721.549999999999
$ws.Range("N122").Value = -22422.2938

# row 133
$ws.Range("H133").Value = 67989
$ws.Range("J133").Value = 67989
$ws.Range("L133").Value = 67989
$ws.Range("N133").Value = -78109

# row 136
$ws.Range("H136").Value = 3794.6
$ws.Range("I136").Value = 2098.75
$ws.Range("J136").Value = 3983.0278
$ws.Range("K136").Value = 6296.25
$ws.Range("L136").Value = 11949.0834
$ws.Range("M136").Value = -3746.25
$ws.Range("N136").Value = -17049.0834
